# Update "paises" (countries) data table and title timestamp
# per upstream data refresh (commit: "Update countries & provincias Spain")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: updated timestamp in title
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 11:04"

# Row 33: Israel
$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 16458
$ws.Range("C33").Value = 4
$ws.Range("D33").Value = 11384
$ws.Range("E33").Value = 4826
$ws.Range("F33").Value = 74
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 248

# Row 35: Polonia
$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 15821
$ws.Range("C35").Value = 170
$ws.Range("D35").Value = 5698
$ws.Range("E35").Value = 9332
$ws.Range("F35").Value = 160
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 791

# Row 36: Japon
$ws.Range("A36").Value = "Japon"
$ws.Range("B36").Value = 15663
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 5906
$ws.Range("E36").Value = 9150
$ws.Range("F36").Value = 287
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 607

# Row 39: Indonesia
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 14032
$ws.Range("C39").Value = 387
$ws.Range("D39").Value = 2698
$ws.Range("E39").Value = 10361
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 14
$ws.Range("H39").Value = 973

# Row 40: Banglades
$ws.Range("A40").Value = "Banglades"
$ws.Range("B40").Value = 13770
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 2414
$ws.Range("E40").Value = 11142
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 214

# Row 42: Filipinas
$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 10794
$ws.Range("C42").Value = 184
$ws.Range("D42").Value = 1924
$ws.Range("E42").Value = 8151
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 15
$ws.Range("H42").Value = 719

# Row 51: Chequia
$ws.Range("A51").Value = "Chequia"
$ws.Range("B51").Value = 8095
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 4448
$ws.Range("E51").Value = 3371
$ws.Range("F51").Value = 40
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 276

# Row 54: Malasia
$ws.Range("A54").Value = "Malasia"
$ws.Range("B54").Value = 6656
$ws.Range("C54").Value = 67
$ws.Range("D54").Value = 5025
$ws.Range("E54").Value = 1523
$ws.Range("F54").Value = 18
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 108

# Row 59: Kazajistan
$ws.Range("A59").Value = "Kazajistan"
$ws.Range("B59").Value = 5056
$ws.Range("C59").Value = 81
$ws.Range("D59").Value = 1828
$ws.Range("E59").Value = 3197
$ws.Range("F59").Value = 31
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 31

# Row 66: Oman
$ws.Range("A66").Value = "Oman"
$ws.Range("B66").Value = 3399
$ws.Range("C66").Value = 175
$ws.Range("D66").Value = 1117
$ws.Range("E66").Value = 2265
$ws.Range("F66").Value = 17
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 17

# Row 67: Armenia
$ws.Range("A67").Value = "Armenia"
$ws.Range("B67").Value = 3313
$ws.Range("C67").Value = 138
$ws.Range("D67").Value = 1325
$ws.Range("E67").Value = 1943
$ws.Range("F67").Value = 10
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 45

# Row 68: Hungria
$ws.Range("A68").Value = "Hungria"
$ws.Range("B68").Value = 3263
$ws.Range("C68").Value = 50
$ws.Range("D68").Value = 933
$ws.Range("E68").Value = 1917
$ws.Range("F68").Value = 50
$ws.Range("G68").Value = 8
$ws.Range("H68").Value = 413

# Row 74: Uzbekistan
$ws.Range("A74").Value = "Uzbekistan"
$ws.Range("B74").Value = 2387
$ws.Range("C74").Value = 38
$ws.Range("D74").Value = 1856
$ws.Range("E74").Value = 521
$ws.Range("F74").Value = 8
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 10

# Row 83: Estonia
$ws.Range("A83").Value = "Estonia"
$ws.Range("B83").Value = 1739
$ws.Range("C83").Value = 6
$ws.Range("D83").Value = 750
$ws.Range("E83").Value = 929
$ws.Range("F83").Value = 5
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 60

# Row 131: Isla de Man
$ws.Range("A131").Value = "Isla de Man"
$ws.Range("B131").Value = 330
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 271
$ws.Range("E131").Value = 36
$ws.Range("F131").Value = 19
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 23

# Row 155: Brunei
$ws.Range("A155").Value = "Brunei"
$ws.Range("B155").Value = 141
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 134
$ws.Range("E155").Value = 6
$ws.Range("F155").Value = 2
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 1

# Row 192: Belice
$ws.Range("A192").Value = "Belice"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 16
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

# Row 193: Nueva Caledonia
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 196: Namibia
$ws.Range("A196").Value = "Namibia"
$ws.Range("B196").Value = 16
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 11
$ws.Range("E196").Value = 5
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0
